# "Updating R Markdown, starting ggplot"
#
# The "Grammar of graphics" / "Themes, labels, facets" module (week 4/5 of
# the "2: Coding fundamentals" unit) is being split: "Grammar of graphics"
# becomes "ggplot 101", and "Themes, labels, facets" becomes
# "Themes, labels, facets (ggplot 102)". Also fixes a typo
# ("2: Coding fundamental" -> "2: Coding fundamentals") on the Schedule
# sheet so it matches the Schedule_date sheet, and nudges the sheet
# selection / column width as the author left them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Schedule_date")

# --- Schedule sheet ---------------------------------------------------

# Typo fix: "2: Coding fundamental" -> "2: Coding fundamentals"
$ws1.Range("B4").Value = "2: Coding fundamentals"
$ws1.Range("B5").Value = "2: Coding fundamentals"
$ws1.Range("B6").Value = "2: Coding fundamentals"

# Split the ggplot module topic names
$ws1.Range("C5").Value = "ggplot 101"
$ws1.Range("C6").Value = "Themes, labels, facets (ggplot 102)"

# --- Schedule_date sheet ----------------------------------------------

$ws2.Range("D5").Value = "ggplot 101"
$ws2.Range("D6").Value = "Themes, labels, facets (ggplot 102)"

# Column C no longer auto-fit; explicit width so "Themes, labels, facets
# (ggplot 102)" isn't clipped.
$ws2.Columns.Item(3).ColumnWidth = 21.666666666666668

# --- Selection / active sheet ------------------------------------------

[void]$ws2.Range("D5:D6").Select()
[void]$ws1.Activate()
[void]$ws1.Range("C5:C6").Select()
